# Update "Training Dashboard" sheet: decrement PERIOD TO EXPIRE (col H) by 1
# and bump LAST UPDATE (col I) from 03-Nov-2025 to 04-Nov-2025 for rows 3-22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Force the LAST UPDATE column to remain plain text so Excel does not
# reinterpret "04-Nov-2025" as a date serial number.
$ws.Range("I3:I22").NumberFormat = "@"

for ($row = 3; $row -le 22; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE

    if ($iCell.Value2 -eq "03-Nov-2025") {
        $hCell.Value2 = $hCell.Value2 - 1
        $iCell.Value2 = "04-Nov-2025"
    }
}
